$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "26.572.95", "1.000") that
# must stay literal text, matching the original inlineStr cells in the workbook.
# Force text format before assigning so Excel does not coerce the string into a
# Double (which would drop significant trailing zeros / reformat the number), then
# restore the cell to the default "Normal" style so no stray per-cell style index
# is left behind (matches the un-styled D/E cells in the source workbook).
function Set-TextValue($cell, $val) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.572.95"
$ws.Range("E2").Value = "  +2.50%  "
Set-TextValue "D3" "1.686.29"
$ws.Range("E3").Value = "  +3.36%  "
$ws.Range("E4").Value = "  -0.44%  "
Set-TextValue "D5" "217.89"
$ws.Range("E5").Value = "  +5.62%  "
Set-TextValue "D6" "0.5364"
$ws.Range("E6").Value = "  +4.93%  "
Set-TextValue "D7" "1.000"
$ws.Range("E7").Value = "  -0.48%  "
Set-TextValue "D8" "0.2685"
$ws.Range("E8").Value = "  +5.78%  "
Set-TextValue "D9" "0.06442"
$ws.Range("E9").Value = "  +4.89%  "
Set-TextValue "D10" "21.38"
$ws.Range("E10").Value = "  +5.39%  "
Set-TextValue "D11" "0.07774"
$ws.Range("E11").Value = "  +2.97%  "
Set-TextValue "D12" "1.685.57"
$ws.Range("E12").Value = "  +0.58%  "
Set-TextValue "D13" "4.501"
$ws.Range("E13").Value = "  +4.01%  "
Set-TextValue "D14" "0.5648"
$ws.Range("E14").Value = "  +6.07%  "
Set-TextValue "D15" "0.0₅8409"
$ws.Range("E15").Value = "  +6.24%  "
Set-TextValue "D16" "66.25"
$ws.Range("E16").Value = "  +1.82%  "
Set-TextValue "D17" "26.589.65"
$ws.Range("E17").Value = "  +2.52%  "
Set-TextValue "D18" "4.833"
$ws.Range("E18").Value = "  +5.26%  "
Set-TextValue "D19" "1.000"
$ws.Range("E19").Value = "  -0.44%  "
Set-TextValue "D20" "195.16"
$ws.Range("E20").Value = "  +5.79%  "
$ws.Range("E21").Value = "  +5.24%  "
Set-TextValue "D22" "6.403"
$ws.Range("E22").Value = "  +6.11%  "
$ws.Range("E23").Value = "  -0.48%  "
Set-TextValue "D24" "143.93"
$ws.Range("E24").Value = "  -1.90%  "
Set-TextValue "D25" "0.1280"
Set-TextValue "D26" "7.508"
$ws.Range("E26").Value = "  +3.50%  "
Set-TextValue "D27" "16.27"
$ws.Range("E27").Value = "  +5.98%  "
Set-TextValue "D28" "1.425"
$ws.Range("E28").Value = "  +5.32%  "
Set-TextValue "D29" "0.06142"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("E30").Value = "  +3.33%  "
Set-TextValue "D31" "3.602"
$ws.Range("E31").Value = "  +6.63%  "
Set-TextValue "D32" "3.469"
$ws.Range("E32").Value = "  +3.92%  "
Set-TextValue "D33" "1.709"
$ws.Range("E33").Value = "  +6.73%  "
Set-TextValue "D34" "1.019"
$ws.Range("E34").Value = "  +6.08%  "
Set-TextValue "D35" "2.800"
$ws.Range("E35").Value = "  +3.15%  "
Set-TextValue "D36" "2.415"
$ws.Range("E36").Value = "  +1.44%  "
Set-TextValue "D37" "0.5749"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  +4.71%  "
Set-TextValue "D39" "5.959"
$ws.Range("E39").Value = "  +3.36%  "
Set-TextValue "D40" "0.8718"
$ws.Range("E40").Value = "  +3.91%  "
Set-TextValue "D41" "1.055.72"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("E42").Value = "  -0.24%  "
Set-TextValue "D43" "100.39"
$ws.Range("E43").Value = "  +0.83%  "
Set-TextValue "D44" "1.836.65"
$ws.Range("E44").Value = "  +2.67%  "
Set-TextValue "D45" "57.35"
$ws.Range("E45").Value = "  +6.44%  "
Set-TextValue "D46" "0.0₈107"
$ws.Range("E46").Value = "  -0.65%  "
Set-TextValue "D47" "8.212"
$ws.Range("E47").Value = "  +4.31%  "
Set-TextValue "D48" "1.002"
$ws.Range("E48").Value = "  +0.23%  "
Set-TextValue "D49" "6.123"
$ws.Range("E49").Value = "  +6.15%  "
Set-TextValue "D50" "0.05204"
$ws.Range("E50").Value = "  +0.08%  "
Set-TextValue "D51" "0.4236"
$ws.Range("E51").Value = "  +0.16%  "
